# "removed duplicate from json file"
#
# C38 and C39 both held the shared string "Placeholder" (the same text
# was duplicated across two FAQ rows). Split them into two distinct,
# non-duplicate strings so each cell carries its own value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C38").Value = "Placeholder 1"
$ws.Range("C39").Value = "Placeholder 2"

# Reflect the author's resulting view/selection state: scrolled down so
# row 6 is at the top, with C40 as the active/selected cell.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C40").Select()
